$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1826923076923077
$ws.Range("C2").Value = 0.5961538461538461
$ws.Range("J2").Value = 0.00641025641025641
$ws.Range("P2").Value = 0.1153846153846154
$ws.Range("S2").Value = 0.09935897435897435
$ws.Range("B3").Value = 0.03125
$ws.Range("C3").Value = 0.03125
$ws.Range("J3").Value = 0.02083333333333333
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("S3").Value = 0.25
$ws.Range("J4").Value = 0.08
$ws.Range("P4").Value = 0.6
$ws.Range("S4").Value = 0.32
$ws.Range("B6").Value = 0.07389162561576355
$ws.Range("D6").Value = 0.01477832512315271
$ws.Range("F6").Value = 0.06403940886699508
$ws.Range("J6").Value = 0.2610837438423645
$ws.Range("O6").Value = 0.01970443349753695
$ws.Range("Q6").Value = 0.1428571428571428
$ws.Range("R6").Value = 0.1330049261083744
$ws.Range("S6").Value = 0.2906403940886699
$ws.Range("B7").Value = 0.1128205128205128
$ws.Range("D7").Value = 0.03076923076923077
$ws.Range("F7").Value = 0.03589743589743589
$ws.Range("J7").Value = 0.1743589743589744
$ws.Range("O7").Value = 0.02051282051282051
$ws.Range("Q7").Value = 0.1948717948717949
$ws.Range("S7").Value = 0.3641025641025641
$ws.Range("B8").Value = 0.1042654028436019
$ws.Range("D8").Value = 0.02132701421800948
$ws.Range("E8").Value = 0.002369668246445498
$ws.Range("F8").Value = 0.05450236966824645
$ws.Range("J8").Value = 0.1540284360189574
$ws.Range("O8").Value = 0.01658767772511848
$ws.Range("Q8").Value = 0.1635071090047393
$ws.Range("R8").Value = 0.07582938388625593
$ws.Range("S8").Value = 0.4075829383886256
$ws.Range("B9").Value = 0.06349206349206349
$ws.Range("D9").Value = 0.01058201058201058
$ws.Range("F9").Value = 0.03174603174603174
$ws.Range("J9").Value = 0.1693121693121693
$ws.Range("O9").Value = 0.02645502645502645
$ws.Range("Q9").Value = 0.164021164021164
$ws.Range("R9").Value = 0.08465608465608465
$ws.Range("S9").Value = 0.4497354497354497
$ws.Range("B10").Value = 0.1050750536097212
$ws.Range("D10").Value = 0.02072909220872051
$ws.Range("E10").Value = 0.0007147962830593281
$ws.Range("F10").Value = 0.0636168691922802
$ws.Range("J10").Value = 0.1486776268763402
$ws.Range("O10").Value = 0.02001429592566119
$ws.Range("Q10").Value = 0.174410293066476
$ws.Range("R10").Value = 0.08434596140100072
$ws.Range("S10").Value = 0.3824160114367405
$ws.Range("G11").Value = 0.124203821656051
$ws.Range("J11").Value = 0.1082802547770701
$ws.Range("K11").Value = 0.2070063694267516
$ws.Range("L11").Value = 0.5445859872611465
$ws.Range("S11").Value = 0.01592356687898089
$ws.Range("G12").Value = 0.7457627118644068
$ws.Range("J12").Value = 0.1807909604519774
$ws.Range("L12").Value = 0.02824858757062147
$ws.Range("S12").Value = 0.04519774011299435
$ws.Range("G13").Value = 0.6222222222222222
$ws.Range("J13").Value = 0.3111111111111111
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("F15").Value = 0.01298701298701299
$ws.Range("H15").Value = 0.09956709956709957
$ws.Range("I15").Value = 0.08658008658008658
$ws.Range("J15").Value = 0.3766233766233766
$ws.Range("K15").Value = 0.03896103896103896
$ws.Range("M15").Value = 0.01298701298701299
$ws.Range("O15").Value = 0.06926406926406926
$ws.Range("S15").Value = 0.303030303030303
$ws.Range("F16").Value = 0.02717391304347826
$ws.Range("H16").Value = 0.1521739130434783
$ws.Range("I16").Value = 0.07065217391304347
$ws.Range("J16").Value = 0.4293478260869565
$ws.Range("K16").Value = 0.108695652173913
$ws.Range("M16").Value = 0.005434782608695652
$ws.Range("N16").Value = 0.005434782608695652
$ws.Range("O16").Value = 0.05978260869565218
$ws.Range("S16").Value = 0.1413043478260869
$ws.Range("F17").Value = 0.007407407407407408
$ws.Range("H17").Value = 0.1802469135802469
$ws.Range("I17").Value = 0.108641975308642
$ws.Range("J17").Value = 0.4074074074074074
$ws.Range("K17").Value = 0.1037037037037037
$ws.Range("M17").Value = 0.009876543209876543
$ws.Range("O17").Value = 0.04197530864197531
$ws.Range("S17").Value = 0.1407407407407407
$ws.Range("F18").Value = 0.02427184466019417
$ws.Range("H18").Value = 0.1601941747572816
$ws.Range("I18").Value = 0.1067961165048544
$ws.Range("J18").Value = 0.4514563106796117
$ws.Range("K18").Value = 0.06796116504854369
$ws.Range("M18").Value = 0.01941747572815534
$ws.Range("N18").Value = 0.004854368932038835
$ws.Range("O18").Value = 0.06796116504854369
$ws.Range("S18").Value = 0.0970873786407767
$ws.Range("F19").Value = 0.01359516616314199
$ws.Range("H19").Value = 0.202416918429003
$ws.Range("I19").Value = 0.06797583081570997
$ws.Range("J19").Value = 0.3829305135951662
$ws.Range("K19").Value = 0.1216012084592145
$ws.Range("M19").Value = 0.02492447129909366
$ws.Range("O19").Value = 0.08081570996978851
$ws.Range("S19").Value = 0.1057401812688822

Write-Output "Applied 108 cell updates"
